$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: (address, new text value)
$updates = @(
    @('D2', '29.837.18'),
    @('E2', '  -0.35%  '),
    @('D3', '1.887.18'),
    @('E3', '  -0.41%  '),
    @('E4', '  +0.03%  '),
    @('D5', '0.7536'),
    @('E5', '  -2.93%  '),
    @('D6', '242.39'),
    @('E6', '  -0.63%  '),
    @('E7', '  -0.04%  '),
    @('D8', '0.3122'),
    @('E8', '  -0.20%  '),
    @('D9', '25.29'),
    @('E9', '  -1.65%  '),
    @('D10', '0.07120'),
    @('E10', '  -3.04%  '),
    @('E11', '  +5.01%  '),
    @('D12', '0.7598'),
    @('E12', '  -1.64%  '),
    @('D13', '1.896.20'),
    @('E13', '  +0.54%  '),
    @('D14', '5.366'),
    @('E14', '  -2.33%  '),
    @('D15', '93.34'),
    @('E15', '  -0.77%  '),
    @('D16', '6.126'),
    @('E16', '  -1.53%  '),
    @('D17', '29.851.22'),
    @('E17', '  -0.23%  '),
    @('D18', '13.70'),
    @('E18', '  -1.97%  '),
    @('D19', '243.35'),
    @('E19', '  -1.60%  '),
    @('D20', '0.000007811'),
    @('E20', '  -0.10%  '),
    @('D21', '0.9991'),
    @('D22', '2.140.69'),
    @('E22', '  +0.39%  '),
    @('D23', '7.998'),
    @('E23', '  -1.59%  '),
    @('D24', '1.001'),
    @('E24', '  +0.10%  '),
    @('E25', '  +0.66%  '),
    @('D26', '9.367'),
    @('E26', '  -1.09%  '),
    @('D27', '162.87'),
    @('E27', '  -0.34%  '),
    @('E28', '  +0.09%  '),
    @('D29', '2.030'),
    @('E29', '  +0.06%  '),
    @('E30', '  +3.56%  '),
    @('D31', '1.533'),
    @('E31', '  -0.88%  '),
    @('D32', '4.505'),
    @('E32', '  +0.62%  '),
    @('D33', '4.163'),
    @('E33', '  +2.48%  '),
    @('D34', '0.05427'),
    @('E34', '  -2.42%  '),
    @('D35', '1.241'),
    @('E35', '  -0.04%  '),
    @('D36', '0.7513'),
    @('E36', '  -0.10%  '),
    @('D37', '1.002'),
    @('E37', '  +0.09%  '),
    @('E38', '  +0.95%  '),
    @('D39', '0.01945'),
    @('E39', '  +0.69%  '),
    @('D40', '2.773'),
    @('E40', '  -0.60%  '),
    @('D41', '0.4463'),
    @('E41', '  -0.14%  '),
    @('D42', '6.102'),
    @('E42', '  +2.27%  '),
    @('D43', '1.097.65'),
    @('E43', '  -0.43%  '),
    @('D44', '72.60'),
    @('E44', '  -1.98%  '),
    @('D45', '0.8615'),
    @('E45', '  +1.25%  '),
    @('D46', '1.000'),
    @('E46', '  -0.05%  '),
    @('D47', '7.713'),
    @('E47', '  +2.47%  '),
    @('D48', '102.32'),
    @('E48', '  +0.14%  '),
    @('D49', '1.859'),
    @('E49', '  -1.46%  '),
    @('D50', '3.038'),
    @('E50', '  +1.51%  '),
    @('D51', '2.036.88'),
    @('E51', '  -0.59%  '),
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $cell = $ws.Range($addr)
    # Force text format so numeric-looking strings (e.g. "1.001") are
    # not auto-converted to actual numbers by Excel, matching the
    # original inline-string (text) cell content.
    $cell.NumberFormat = "@"
    $cell.Value = $val
}
